# Applies updated profit-calculation figures (currentAveragePrice / LevePrice / LeveProfit
# columns H..N) across all 8 job sheets, per the scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 604.7143
$ws.Range("J17").Value = 612.07275
$ws.Range("L17").Value = 1836.21825
$ws.Range("N17").Value = -2172.21825
$ws.Range("H112").Value = 2598.2144
$ws.Range("J112").Value = 2732.6924
$ws.Range("L112").Value = 8198.0772
$ws.Range("N112").Value = -10414.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 297048.75
$ws.Range("I132").Value = 457038.47
$ws.Range("J132").Value = 3734.3333
$ws.Range("K132").Value = 1371115.41
$ws.Range("L132").Value = 11202.9999
$ws.Range("M132").Value = -1368585.41
$ws.Range("N132").Value = -16262.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13960
$ws.Range("I26").Value = 8440
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 8440
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = -8148
$ws.Range("N26").Value = -25584
$ws.Range("H86").Value = 1685.6842
$ws.Range("I86").Value = 1775.5862
$ws.Range("J86").Value = 1396
$ws.Range("K86").Value = 1775.5862
$ws.Range("L86").Value = 1396
$ws.Range("M86").Value = -652.5862
$ws.Range("N86").Value = -3642
$ws.Range("H89").Value = 1685.6842
$ws.Range("I89").Value = 1775.5862
$ws.Range("J89").Value = 1396
$ws.Range("K89").Value = 8877.931
$ws.Range("L89").Value = 6980
$ws.Range("M89").Value = -3261.931
$ws.Range("N89").Value = -18212
$ws.Range("H96").Value = 4328.5
$ws.Range("I96").Value = 4328.5
$ws.Range("K96").Value = 4328.5
$ws.Range("M96").Value = -1582.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1633.0416
$ws.Range("I122").Value = 1113.2106
$ws.Range("J122").Value = 3608.4
$ws.Range("K122").Value = 3339.6318
$ws.Range("L122").Value = 10825.2
$ws.Range("M122").Value = -889.6318000000001
$ws.Range("N122").Value = -15725.2
$ws.Range("H132").Value = 3989.2222
$ws.Range("I132").Value = 3738
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 11214
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -8684
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1266.6666
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 4500
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -4838
$ws.Range("H27").Value = 1266.6666
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 4500
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -4704
$ws.Range("H32").Value = 463.75
$ws.Range("J32").Value = 463.75
$ws.Range("L32").Value = 1391.25
$ws.Range("N32").Value = -1957.25
$ws.Range("H39").Value = 2491.7144
$ws.Range("J39").Value = 2491.7144
$ws.Range("L39").Value = 7475.1432
$ws.Range("N39").Value = -8063.1432
$ws.Range("H46").Value = 2833.3333
$ws.Range("J46").Value = 7000
$ws.Range("L46").Value = 21000
$ws.Range("N46").Value = -21182
$ws.Range("H58").Value = 2581.6667
$ws.Range("I58").Value = 2472.5
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 7417.5
$ws.Range("L58").Value = 8400
$ws.Range("M58").Value = -7289.5
$ws.Range("N58").Value = -8656
$ws.Range("H64").Value = 7435.778
$ws.Range("I64").Value = 637.3333
$ws.Range("J64").Value = 8795.467000000001
$ws.Range("K64").Value = 1911.9999
$ws.Range("L64").Value = 26386.401
$ws.Range("M64").Value = -1641.9999
$ws.Range("N64").Value = -26926.401
$ws.Range("H67").Value = 7435.778
$ws.Range("I67").Value = 637.3333
$ws.Range("J67").Value = 8795.467000000001
$ws.Range("K67").Value = 1911.9999
$ws.Range("L67").Value = 26386.401
$ws.Range("M67").Value = -975.9999
$ws.Range("N67").Value = -28258.401
$ws.Range("H70").Value = 5159.7896
$ws.Range("I70").Value = 2620.25
$ws.Range("J70").Value = 7006.727
$ws.Range("K70").Value = 7860.75
$ws.Range("L70").Value = 21020.181
$ws.Range("M70").Value = -7545.75
$ws.Range("N70").Value = -21650.181
$ws.Range("H73").Value = 5159.7896
$ws.Range("I73").Value = 2620.25
$ws.Range("J73").Value = 7006.727
$ws.Range("K73").Value = 7860.75
$ws.Range("L73").Value = 21020.181
$ws.Range("M73").Value = -6768.75
$ws.Range("N73").Value = -23204.181
$ws.Range("H76").Value = 7205.7646
$ws.Range("J76").Value = 8035.5713
$ws.Range("L76").Value = 24106.7139
$ws.Range("N76").Value = -24872.7139
$ws.Range("H79").Value = 7205.7646
$ws.Range("J79").Value = 8035.5713
$ws.Range("L79").Value = 24106.7139
$ws.Range("N79").Value = -26758.7139
$ws.Range("H100").Value = 7663.0557
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 8245.9375
$ws.Range("K100").Value = 9000
$ws.Range("L100").Value = 24737.8125
$ws.Range("M100").Value = -8189
$ws.Range("N100").Value = -26359.8125
$ws.Range("H106").Value = 6800
$ws.Range("J106").Value = 6800
$ws.Range("L106").Value = 20400
$ws.Range("N106").Value = -22292
$ws.Range("H123").Value = 5190
$ws.Range("J123").Value = 6400
$ws.Range("L123").Value = 19200
$ws.Range("N123").Value = -24100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3087.375
$ws.Range("I80").Value = 2922
$ws.Range("J80").Value = 3227.3076
$ws.Range("K80").Value = 2922
$ws.Range("L80").Value = 3227.3076
$ws.Range("M80").Value = -1924
$ws.Range("N80").Value = -5223.3076
$ws.Range("H83").Value = 3087.375
$ws.Range("I83").Value = 2922
$ws.Range("J83").Value = 3227.3076
$ws.Range("K83").Value = 14610
$ws.Range("L83").Value = 16136.538
$ws.Range("M83").Value = -9618
$ws.Range("N83").Value = -26120.538
$ws.Range("H132").Value = 2397.4
$ws.Range("I132").Value = 2033.9375
$ws.Range("J132").Value = 3043.5557
$ws.Range("K132").Value = 6101.8125
$ws.Range("L132").Value = 9130.667099999999
$ws.Range("M132").Value = -3571.8125
$ws.Range("N132").Value = -14190.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3883
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 4824.5
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 14473.5
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -19533.5
$ws.Range("H133").Value = 46494.75
$ws.Range("J133").Value = 46494.75
$ws.Range("L133").Value = 46494.75
$ws.Range("N133").Value = -51554.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4235
$ws.Range("I62").Value = 4445
$ws.Range("J62").Value = 4165
$ws.Range("K62").Value = 4445
$ws.Range("L62").Value = 4165
$ws.Range("M62").Value = -3821
$ws.Range("N62").Value = -5413
$ws.Range("H65").Value = 4235
$ws.Range("I65").Value = 4445
$ws.Range("J65").Value = 4165
$ws.Range("K65").Value = 22225
$ws.Range("L65").Value = 20825
$ws.Range("M65").Value = -19105
$ws.Range("N65").Value = -27065
$ws.Range("H81").Value = 7439.8823
$ws.Range("I81").Value = 20755.2
$ws.Range("J81").Value = 1891.8334
$ws.Range("K81").Value = 41510.4
$ws.Range("L81").Value = 3783.6668
$ws.Range("M81").Value = -40449.4
$ws.Range("N81").Value = -5905.6668
$ws.Range("H84").Value = 7439.8823
$ws.Range("I84").Value = 20755.2
$ws.Range("J84").Value = 1891.8334
$ws.Range("K84").Value = 207552
$ws.Range("L84").Value = 18918.334
$ws.Range("M84").Value = -202248
$ws.Range("N84").Value = -29526.334
